# Write results and discussion
# Adds a "DIFF" comparison block (columns W:AF) to the "OPN" sheet that
# mirrors the existing WT (A:J) / QTY (L:U) staircase comparison tables,
# computing QTY-minus-WT differences wherever both sides have numeric data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)   # "OPN" sheet

# --- Pre-seed cell formatting for the new block by cloning the existing
#     "QTY" staircase block (L1:U10) styles onto the new range (W1:AF10).
$ws.Range("L1:U10").Copy()
$ws.Range("W1:AF10").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("L1").Copy()
$ws.Range("W1").PasteSpecial(-4122)        # header-label style (same as L1/A1)

# --- Header row (row 1) ---
$ws.Range("W1").Value = "DIFF"
$ws.Range("X1").Value = "OPN1MW"
$ws.Range("Y1").Value = "OPN1LW"
$ws.Range("Z1").Value = "OPN1SW"
$ws.Range("AA1").Value = "OPN2"
$ws.Range("AB1").Value = "OPN3"
$ws.Range("AC1").Value = "OPN4"
$ws.Range("AD1").Value = "OPN5"
$ws.Range("AE1").Value = "RGR"
$ws.Range("AF1").Value = "RRH"

# --- Row labels (column W, rows 2-10) - same labels as columns A / L ---
$labels = @("OPN1MW","OPN1LW","OPN1SW","OPN2","OPN3","OPN4","OPN5","RGR","RRH")
for ($i = 0; $i -lt 9; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r, 23).Value = $labels[$i]   # column 23 = W
}

# --- Marker column (X, rows 2-10) mirrors columns B / M: always "-" ---
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 24).Value = "-"           # column 24 = X
}

# --- Staircase of differences (Y:AF, rows 2-10) ---
# Destination columns Y..AF line up with source columns C..J (WT block)
# and N..U (QTY block) at the same relative position.
$destCols  = @("Y","Z","AA","AB","AC","AD","AE","AF")
$leftCols  = @("C","D","E","F","G","H","I","J")
$rightCols = @("N","O","P","Q","R","S","T","U")

for ($r = 2; $r -le 10; $r++) {
    for ($i = 0; $i -lt 8; $i++) {
        $dcell = $destCols[$i] + $r
        if ($i -lt ($r - 2)) {
            # Below/on the diagonal: source cells are "-" placeholders too.
            $ws.Range($dcell).Value = "-"
        } else {
            $formula = "=" + $rightCols[$i] + $r + "-" + $leftCols[$i] + $r
            $ws.Range($dcell).Formula = $formula
        }
    }
}

# --- Sheet view / window cosmetics to match the saved state ---
$ws.Range("Y1:Y10").Columns.AutoFit()
$excel.ActiveWindow.Zoom = 55
$ws.Range("AH18").Select()
